# Update "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets (regenerated report).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 3 (0f3cfe74-...) and Row 4 (60069378-...) shared the same
# Handoff/Handback timestamps before, and still share the same
# (refreshed) timestamps now.
$wsZh.Range("E3").Value = "2016-03-12 10:15:47"
$wsZh.Range("E4").Value = "2016-03-12 10:15:47"
$wsZh.Range("H3").Value = "2016-03-12 10:16:10"
$wsZh.Range("H4").Value = "2016-03-12 10:16:10"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("E3").Value = "2016-03-12 10:15:51"
$wsDe.Range("E4").Value = "2016-03-12 10:15:51"
$wsDe.Range("H3").Value = "2016-03-12 10:16:16"
$wsDe.Range("H4").Value = "2016-03-12 10:16:16"
